$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows where "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
$dnasrRows = @(3,6,10,11,12,13,14,15,17,18,30,33,37,38,39,40,41,42,44,45,57,60,64,65,66,67,68,69,71,72,86,87,88,89,93,95,96,112,113,114,115,119,121,122,138,139,140,141,145,147,148)
foreach ($r in $dnasrRows) {
    $ws.Cells.Item($r, 7).Value = "dnasr281@gmail.com, System"
}

# Rows where "System, backup@backdoor.com, system" -> "system, System, backup@backdoor.com"
$backdoorRows = @(2,29,56)
foreach ($r in $backdoorRows) {
    $ws.Cells.Item($r, 7).Value = "system, System, backup@backdoor.com"
}
